# Update "Test Cases" worksheet to reflect the completed test run:
#  - retitle the suite to mention "json" data
#  - rename the function under test
#  - bump the "last modified" timestamp
#  - point the test step at the real test fixture file
#  - flip both test-step statuses from FAILED to PASSED
#  - adjust the active view (zoom + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SS-read-write-data-json")

# Title banner
$ws.Range("A1").Value = "Test Suite 1 - Automated testing of functions read/write json data"

# Test Details block
$ws.Range("D5").Value = "getStatisticsFrom1"
$ws.Range("F7").Value = "08/04/2021 02:41 PM"

# Test steps table
$ws.Range("E11").Value = "Test.json"
$ws.Range("E11").WrapText = $true
$ws.Range("H11").Value = "PASSED"
$ws.Range("H12").Value = "PASSED"

# View state: zoom + active cell selection
$ws.Select()
$excel.ActiveWindow.Zoom = 115
$ws.Range("C22").Select()
